# Add test case "RNCR1" (resolve new card request, print a new card) to the
# "2.9. <Staff> Resolve new card request" section.
#
# In the original sheet, row 17 and row 18 are both blank spacer rows
# (between the "2.8. <Staff> View new card requests" table and the
# "2.9. <Staff> Resolve new card request" header on row 19). The edit
# removes the extra blank spacer row (row 18) - shifting every following
# row up by one - and fills the remaining blank row (row 17) with the new
# test-case data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-redundant second blank spacer row; this shifts rows
# 19..999 up to 18..998 (dimension shrinks from Z999 to Z998).
$ws.Rows.Item(18).Delete()

# Populate the remaining blank row 17 with the new RNCR1 test case.
# Columns: A=ID, B=Test Case Description, C=Test case procedure,
# D=Expected output, E=Inter-test case dependence, F=Result, G=Test Date.
# (E and D are written in this order so the new shared strings land at the
# same indices as the authored workbook: 83..87.)
$ws.Range("A17").Value = "RNCR1"
$ws.Range("B17").Value = "Test resolve new card request, print a new card"
$ws.Range("C17").Value = "1. Contract HD0001 is Ready and having an actived card.`n2. Create new card request for contract HD001 like test case SNCR1.`n3. Print card for Ready contract ""HD0001"" like test case PC2."
$ws.Range("E17").Value = "SNCR1,PC2."
$ws.Range("D17").Value = "Printer app show message ""Write success"", contract card change to new card ID, old card is deactivated"
$ws.Range("F17").Value = "Pass"

# Test Date column uses the same date format as the rest of the sheet
# (copy the format from an existing date cell so it reuses the same style
# instead of Excel minting a brand-new number format).
$ws.Range("G3").Copy()
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G17").Value = 42198

# Keep the row at its original (non-autofit) height.
$ws.Rows.Item(17).RowHeight = 15.75

# Restore the view: scrolled so row 13 is at the top, with G17 selected.
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("G17").Select()
